# Update the average_county_temperature column (AA) with refreshed NOAA
# temperature data for the affected facilities. Each contiguous block of
# rows below corresponds to a single facility_id, which shares a single
# updated temperature value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Range = "AA12:AA16";   Value = 12.51681286549706 },
    @{ Range = "AA17:AA21";   Value = 15.74228395061728 },
    @{ Range = "AA22:AA61";   Value = 1.925925925925943 },
    @{ Range = "AA62:AA71";   Value = 12.66820987654322 },
    @{ Range = "AA77:AA81";   Value = -3.222222222222223 },
    @{ Range = "AA114:AA118"; Value = 1.925925925925943 },
    @{ Range = "AA119:AA128"; Value = 20.68981481481483 },
    @{ Range = "AA129:AA138"; Value = 14.96875 },
    @{ Range = "AA139:AA163"; Value = 1.925925925925943 },
    @{ Range = "AA164:AA173"; Value = -3.222222222222223 }
)

foreach ($u in $updates) {
    $ws.Range($u.Range).Value = $u.Value
}
